$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = 0.0722
$ws.Range("G2").Value = -0.005844027640671274
$ws.Range("H2").Value = -0.005844027640671274
$ws.Range("I2").Value = -0.05478775913129319
$ws.Range("J2").Value = -0.05478775913129319
$ws.Range("K2").Value = -2.73
$ws.Range("L2").Value = -0.02694965449160908
$ws.Range("M2").Value = 0
$ws.Range("N2").Value = 0
$ws.Range("O2").Value = 0
$ws.Range("P2").Value = 0
$ws.Range("Q2").Value = 0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = 0
$ws.Range("T2").ClearContents()
$ws.Range("U2").Value = 15
$ws.Range("V2").Value = 0.07804370447450573
$ws.Range("W2").Value = -0.0303670745272525
$ws.Range("X2").Value = 0.1001083537064584
$ws.Range("Y2").Value = -0.1304754282337109
$ws.Range("Z2").Value = 0.4172158154859967
$ws.Range("AA2").Value = -0.02285831960461285
$ws.Range("AB2").Value = 0.06503686196336898
$ws.Range("AC2").Value = -0.08789518156798183
$ws.Range("AD2").Value = 203.6
$ws.Range("AE2").Value = 0
$ws.Range("AF2").Value = 203.6
$ws.Range("AG2").Value = 188.6
$ws.Range("AH2").Value = 0.514401212733704
$ws.Range("AI2").Value = 0.7474302496328928
$ws.Range("AJ2").Value = 0.4952731092436975
$ws.Range("AK2").Value = 0.7327117327117327
$ws.Range("AL2").Value = 1.97
$ws.Range("AM2").Value = 1.796
$ws.Range("AN2").Value = 43.59743040685225
$ws.Range("AO2").Value = -2.817258883248731
$ws.Range("AP2").Value = 40.38543897216274
$ws.Range("AQ2").Value = -3.090200445434298

# Row 3
$ws.Range("D3").Value = 0.0722
$ws.Range("G3").Value = -0.005844027640671274
$ws.Range("H3").Value = -0.005844027640671274
$ws.Range("I3").Value = -0.05478775913129319
$ws.Range("J3").Value = -0.05478775913129319
$ws.Range("K3").Value = -2.73
$ws.Range("L3").Value = -0.02694965449160908
$ws.Range("M3").Value = 0
$ws.Range("N3").Value = 0
$ws.Range("O3").Value = 0
$ws.Range("P3").Value = 0
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").ClearContents()
$ws.Range("U3").Value = 15
$ws.Range("V3").Value = 0.07804370447450573
$ws.Range("W3").Value = -0.0303670745272525
$ws.Range("X3").Value = 0.1001083537064584
$ws.Range("Y3").Value = -0.1304754282337109
$ws.Range("Z3").Value = 0.4172158154859967
$ws.Range("AA3").Value = -0.02285831960461285
$ws.Range("AB3").Value = 0.06503686196336898
$ws.Range("AC3").Value = -0.08789518156798183
$ws.Range("AD3").Value = 203.6
$ws.Range("AE3").Value = 0
$ws.Range("AF3").Value = 203.6
$ws.Range("AG3").Value = 188.6
$ws.Range("AH3").Value = 0.514401212733704
$ws.Range("AI3").Value = 0.7474302496328928
$ws.Range("AJ3").Value = 0.4952731092436975
$ws.Range("AK3").Value = 0.7327117327117327
$ws.Range("AL3").Value = 1.97
$ws.Range("AM3").Value = 1.796
$ws.Range("AN3").Value = 43.59743040685225
$ws.Range("AO3").Value = -2.817258883248731
$ws.Range("AP3").Value = 40.38543897216274
$ws.Range("AQ3").Value = -3.090200445434298
